$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 99, shifting rows 99:149 down to 100:150
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with data (matching the template of surrounding rows)
$ws.Cells.Item(99, 1).Value = 9
$ws.Cells.Item(99, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(99, 3).Value = "Metropolitana"
$ws.Cells.Item(99, 4).Value = 44455
$ws.Cells.Item(99, 5).Value = 13
$ws.Cells.Item(99, 6).Value = 300000001
$ws.Cells.Item(99, 7).Value = "Rabanito"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 7900
$ws.Cells.Item(99, 11).Value = 3500
$ws.Cells.Item(99, 12).Value = 4000
$ws.Cells.Item(99, 13).Value = 3747
$ws.Cells.Item(99, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(99, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(99, 16).Value = 37
$ws.Cells.Item(99, 17).Value = 100
$ws.Cells.Item(99, 18).Value = "Hortaliza"
